$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.139.40'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '2.731.57'
$ws.Range('E3').Value = '  -5.62%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'507.64"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.38%  '
$ws.Range('D6').Value = "'141.49"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').Value = "'0.534"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.64%  '
$ws.Range('D9').Value = '2.747.00'
$ws.Range('E9').Value = '  -4.98%  '
$ws.Range('D10').Value = "'6.13"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('D12').Value = "'0.349"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').Value = '3.211.34'
$ws.Range('E14').Value = '  -5.51%  '
$ws.Range('D15').Value = '59.037.92'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('D16').Value = "'21.85"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.28%  '
$ws.Range('D17').Value = "'0.0000136"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '2.739.09'
$ws.Range('E18').Value = '  -5.33%  '
$ws.Range('D19').Value = "'4.75"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('D20').Value = "'11.01"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.10%  '
$ws.Range('D21').Value = "'346.13"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.81%  '
$ws.Range('D22').Value = "'6.26"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.42%  '
$ws.Range('D23').Value = "'0.998"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('D25').Value = "'63.28"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = "'0.174"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').Value = "'0.425"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.04%  '
$ws.Range('D28').Value = "'0.990"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('D29').Value = '0.0₃0840'
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('D30').Value = "'7.50"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.09%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').Value = "'1.63"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').Value = "'19.15"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').Value = "'149.08"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('D35').Value = "'4.22"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('D36').Value = "'5.40"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.11%  '
$ws.Range('D37').Value = "'0.951"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.84%  '
$ws.Range('D38').Value = "'1.14"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('D39').Value = "'36.21"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.62%  '
$ws.Range('D40').Value = "'1.39"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.95%  '
$ws.Range('D41').Value = "'3.54"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.91%  '
$ws.Range('D42').Value = '2.190.76'
$ws.Range('E42').Value = '  -6.11%  '
$ws.Range('D43').Value = "'0.0559"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').Value = "'0.996"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = "'0.606"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.40%  '
$ws.Range('D46').Value = "'19.16"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.21%  '
$ws.Range('D47').Value = "'4.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D49').Value = "'0.0229"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('D50').Value = "'0.0887"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.19%  '
$ws.Range('D51').Value = "'18.09"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.22%  '
